# "replaced (A) and reloaded"
#
# The "Staff Performance Overview" merged report carried separate "(A)"
# adjustment rows for two staff members (Alicia Gallegos (A), Brook
# Accomando (A)). This reload folds each "(A)" row's totals back into the
# staff member's primary row and removes the now-redundant "(A)" rows, which
# pulls every following row up by one or two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 13   # column M - last numeric column that gets summed (N/O are derived averages)
$lastRowBefore = 22

# Snapshot every data row (3..22) from the sheet *before* any writes happen,
# so later writes never clobber a value we still need to read.
$snapshot = @{}
for ($r = 3; $r -le $lastRowBefore; $r++) {
    $row = @{}
    $row[1] = $ws.Cells.Item($r, 1).Value2
    for ($c = 2; $c -le $lastCol; $c++) {
        $row[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $row
}

# Destination row (in the reloaded sheet) -> source row(s) in the original
# sheet whose values are combined into it. Rows with two sources are the
# staff member's row plus its "(A)" adjustment row, folded together.
$sourceMap = @{
    3  = @(3)
    4  = @(4, 5)
    5  = @(6)
    6  = @(7, 8)
    7  = @(9)
    8  = @(10)
    9  = @(11)
    10 = @(12)
    11 = @(13)
    12 = @(14)
    13 = @(15)
    14 = @(16)
    15 = @(17)
    16 = @(18)
    17 = @(19)
    18 = @(20)
    19 = @(21)
    20 = @(22)
}

foreach ($destRow in ($sourceMap.Keys | Sort-Object)) {
    $sources = $sourceMap[$destRow]

    # Name comes from the first (primary) source row.
    $ws.Cells.Item($destRow, 1).Value2 = $snapshot[$sources[0]][1]

    # Sum each numeric column (B..M) across the source row(s).
    for ($c = 2; $c -le $lastCol; $c++) {
        $total = 0
        foreach ($src in $sources) {
            $total += $snapshot[$src][$c]
        }
        $ws.Cells.Item($destRow, $c).Value2 = $total
    }

    # Recompute the derived averages: Avg. Spend Per Client = Total / Client Visits #
    $clients = $ws.Cells.Item($destRow, 2).Value2
    $ws.Cells.Item($destRow, 14).Value2 = $ws.Cells.Item($destRow, 12).Value2 / $clients
    $ws.Cells.Item($destRow, 15).Value2 = $ws.Cells.Item($destRow, 13).Value2 / $clients
}

# The reload is two rows shorter than the original (the two "(A)" rows are
# gone) - clear the now-unused trailing rows so the sheet's used range and
# dimension shrink to match.
$ws.Range("A21:O22").Clear()
